$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19089138507843
$ws.Range("B1").Value = 2.474459648132324
$ws.Range("D1").Value = 2.277533531188965
$ws.Range("E1").Value = 1.179970979690552
